# Actualización automática desde WSL
# Appends the latest batch of DropControl sensor readings (week/"Sem" 32,
# 2025-06-08 18:00 through 23:00) to Sheet1 and refreshes the existing
# week-number column plus the active selection / print orientation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) The whole existing data set belongs to ISO week 32 now (previously
#    stamped as 23) -- update column C ("Sem") for the existing rows.
# ---------------------------------------------------------------------
for ($r = 2; $r -le 14; $r++) {
    $ws.Cells.Item($r, 3).Value = 32
}

# ---------------------------------------------------------------------
# 2) Append six new hourly readings (rows 15-20) below the current data.
#    Column A keeps the same custom date-time display format used by the
#    rows above it (A2:A14), column H stores the predominant wind
#    direction as text ("E" / "ESE", already present in the shared
#    string table).
# ---------------------------------------------------------------------
$ws.Range("A15:A20").NumberFormat = "yyyy\-mm\-dd\ hh:mm:ss"

$newRows = @(
    @{ Row=15; A=45816.75;              B=2025; C=32; D=16.43; E=84.13; F=0.88; G=4.61; H="E"   ; I=0 },
    @{ Row=16; A=45816.791666666664;    B=2025; C=32; D=15.04; E=87.73; F=0;    G=3.93; H="E"   ; I=0 },
    @{ Row=17; A=45816.833333333336;    B=2025; C=32; D=14.39; E=90.48; F=0;    G=4.05; H="ESE" ; I=0 },
    @{ Row=18; A=45816.875;             B=2025; C=32; D=14.35; E=91.32; F=0;    G=0.5;  H="E"   ; I=0 },
    @{ Row=19; A=45816.916666666664;    B=2025; C=32; D=14.32; E=91.42; F=0;    G=1.65; H="E"   ; I=0 },
    @{ Row=20; A=45816.958333333336;    B=2025; C=32; D=14.14; E=91.56; F=0;    G=0.73; H="ESE" ; I=0 }
)

foreach ($row in $newRows) {
    $r = $row.Row
    $ws.Cells.Item($r, 1).Value = $row.A
    $ws.Cells.Item($r, 2).Value = $row.B
    $ws.Cells.Item($r, 3).Value = $row.C
    $ws.Cells.Item($r, 4).Value = $row.D
    $ws.Cells.Item($r, 5).Value = $row.E
    $ws.Cells.Item($r, 6).Value = $row.F
    $ws.Cells.Item($r, 7).Value = $row.G
    $ws.Cells.Item($r, 8).Value = $row.H
    $ws.Cells.Item($r, 9).Value = $row.I
}

# ---------------------------------------------------------------------
# 3) Restore the print orientation (portrait) and move the active cell
#    selection to reflect where the operator's cursor ended up.
# ---------------------------------------------------------------------
$ws.PageSetup.Orientation = 1

$ws.Range("F13").Select()

# Best-effort: reflect the window position recorded by Excel on save.
$wb.Windows.Item(1).Left = -105
